# Update the login/credentials demo data on Sheet1 and move the selection,
# matching the "added demo test files" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new manager id / generated password
$ws.Range("A2").Value = "mngr429679"
$ws.Range("B2").Value = "jYdyvYg"

# Row 3: same new manager id repeated in both columns
$ws.Range("A3").Value = "mngr429679"
$ws.Range("B3").Value = "mngr429679"

# Row 4: username stays the same, password replaced
$ws.Range("B4").Value = "rahjjjjjaguh"

# Leave the current selection on B3, as last left by the editor
$ws.Range("B3").Select()
